$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.813.81"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.089.35"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.67"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.39"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.388"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.20%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.385.05"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.45"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.11"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.46%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.24"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.100.80"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.701.38"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.15"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.94"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.20"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.83%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.47"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.141"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +11.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.94"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.44"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0625"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.46%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.50"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.38"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.80%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.69%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.11"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.88%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.454.19"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.10"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.31%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.06"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.64"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.66%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.23%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.280.41"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.08%  "
